$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'isusana481'
$ws.Range("B3").Value = 'lelena520'
$ws.Range("B4").Value = 'jvargas250'
$ws.Range("B5").Value = 'fchacón520'
$ws.Range("B6").Value = 'avargas546'
$ws.Range("B7").Value = 'arojas606'
$ws.Range("B8").Value = 'cvargas910'
$ws.Range("B9").Value = 'halpízar088'
$ws.Range("B10").Value = 'rálvarez670'
$ws.Range("B11").Value = 'scamacho014'
$ws.Range("B12").Value = 'hvíquez221'
$ws.Range("B13").Value = 'ppérez191'
$ws.Range("B14").Value = 'augalde080'
$ws.Range("B15").Value = 'ccortés020'
$ws.Range("B16").Value = 'lcastillo300'
$ws.Range("B17").Value = 'cdíaz294'
$ws.Range("B18").Value = 'yhernández601'
$ws.Range("B19").Value = 'vsantos171'
$ws.Range("B20").Value = 'gbarrantes885'
$ws.Range("B21").Value = 'dálvarez007'
$ws.Range("B22").Value = 'vsegura020'
$ws.Range("B23").Value = 'pcalvo201'
$ws.Range("B24").Value = 'yvargas925'
$ws.Range("B25").Value = 'gsolano301'
$ws.Range("B26").Value = 'omarín460'
$ws.Range("B27").Value = 'yabarca440'
$ws.Range("B28").Value = 'mfallas490'
$ws.Range("B29").Value = 'schavarría522'
$ws.Range("B30").Value = 'asánchez641'
$ws.Range("B31").Value = 'mzamora622'
$ws.Range("B32").Value = 'abermúdez019'
$ws.Range("B33").Value = 'yhernández100'
$ws.Range("B34").Value = 'lsalas339'
$ws.Range("B35").Value = 'harce041'
$ws.Range("B36").Value = 'marguello265'
$ws.Range("B37").Value = 'serasmo074'
$ws.Range("B38").Value = 'agonzález230'
$ws.Range("B39").Value = 'jmaría795'
$ws.Range("B40").Value = 'cramírez181'
$ws.Range("B41").Value = 'kbolaños340'
$ws.Range("B42").Value = 'cprado320'
$ws.Range("B44").Value = 'fsánchez400'
$ws.Range("B45").Value = 'kvargas870'
$ws.Range("B46").Value = 'mvalverde872'
$ws.Range("B47").Value = 'kcampos777'
$ws.Range("B48").Value = 'ycórdoba257'
$ws.Range("B49").Value = 'cyesenia127'
$ws.Range("B50").Value = 'ymaría614'
$ws.Range("B51").Value = 'osegura608'
$ws.Range("B52").Value = 'rsoto066'
$ws.Range("B53").Value = 'jcontreras061'
$ws.Range("B54").Value = 'macuña000'
$ws.Range("B55").Value = 'isalas691'
$ws.Range("B57").Value = 'kbonilla032'
$ws.Range("B58").Value = 'gquesada753'
$ws.Range("B59").Value = 'mcatalina033'
$ws.Range("B60").Value = 'jrodríguez162'
$ws.Range("B61").Value = 'yruiz538'
$ws.Range("B62").Value = 'rsalas240'
$ws.Range("B63").Value = 'ecambronero404'
$ws.Range("B64").Value = 'mdenise675'
$ws.Range("B65").Value = 'jsalazar227'
$ws.Range("B66").Value = 'alorena066'
$ws.Range("B67").Value = 'frojas871'
$ws.Range("B68").Value = 'karguedas689'
$ws.Range("B69").Value = 'asalas328'
$ws.Range("B70").Value = 'bsteven562'
$ws.Range("B71").Value = 'jarguedas047'
$ws.Range("B72").Value = 'despinoza095'
$ws.Range("B73").Value = 'cquirós062'
$ws.Range("B74").Value = 'jarturo773'
$ws.Range("B76").Value = 'fgonzález715'
$ws.Range("B77").Value = 'mgonzález080'
$ws.Range("B78").Value = 'hoviedo895'
$ws.Range("B79").Value = 'ycastillo200'
$ws.Range("B80").Value = 'kmiranda410'
$ws.Range("B81").Value = 'apalma929'
$ws.Range("B82").Value = 'rrojas576'
$ws.Range("B83").Value = 'jazofeifa047'
$ws.Range("B84").Value = 'pvargas778'
$ws.Range("B85").Value = 'adíaz004'
$ws.Range("B86").Value = 'hgamboa201'
$ws.Range("B87").Value = 'csolís107'
$ws.Range("B88").Value = 'jmurillo305'
$ws.Range("B89").Value = 'amaría005'
$ws.Range("B90").Value = 'equesada868'
$ws.Range("B91").Value = 'ggonzález302'
$ws.Range("B92").Value = 'jmanuel893'
$ws.Range("B93").Value = 'rmurillo470'
$ws.Range("B94").Value = 'kbryan084'
$ws.Range("B95").Value = 'gbrenes573'
$ws.Range("B96").Value = 'gsoto125'
$ws.Range("B97").Value = 'wrodríguez500'
$ws.Range("B98").Value = 'mramírez170'
$ws.Range("B99").Value = 'erodríguez904'
$ws.Range("B100").Value = 'koviedo306'
$ws.Range("B101").Value = 'dvega363'
$ws.Range("B102").Value = 'abonilla166'
$ws.Range("B103").Value = 'kjaen199'
$ws.Range("B104").Value = 'mmagdalena905'
